$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new row 12: phone "79174414" stored as text (matches source data
# where this record's phone wasn't parsed as numeric), birthday left blank,
# total_points set to 0.
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = "79174414"
$ws.Range("B12").Value = ""
$ws.Range("C12").Value = 0
